$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Write the date as literal text into A6. Excel (and this engine) would
# normally auto-convert a bare "2024-10-05" into a date serial + a new
# number-format style, which the target workbook does not have (it stores
# plain shared-string text, like the other Date cells A2:A5). To avoid the
# autoconversion we enter the value with a trailing space (never matches a
# date pattern, so it stays text), then TRIM it via a helper formula cell
# and paste back just the resulting value - that round trip keeps the cell
# a plain text value with no special number format applied.
$ws.Range("A6").Value = "2024-10-05 "

$helper = $ws.Range("D1")
$helper.Formula = "=TRIM(A6)"
$helper.Copy()
$ws.Range("A6").PasteSpecial(-4163)  # xlPasteValues
$helper.ClearContents()
$excel.CutCopyMode = 0

$ws.Range("B6").Value = 0.9983
